# Applies the "Updated symbol list ... with GitHub Actions" refresh to the
# cryptos worksheet: refreshed Price values for most rows, plus rows 18-24
# whose Coin/Link/Price/Volume(1h) columns rolled over to the next ranking
# entry.
#
# All of the touched cells in the workbook are stored as text (the sheet
# uses inline/shared strings everywhere, even for the "Price" column), so a
# plain numeric-looking assignment such as $rng.Value = '248.76' must be
# protected from Excel's automatic number conversion. Set-TextValue forces
# the cell to Text format before assigning the string, then restores the
# cell's (default) style so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$Address, [string]$Val)
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

# --- Price-only refreshes (rows 2-17) ---
Set-TextValue "D2"  "248.76"
Set-TextValue "D3"  "21.74"
Set-TextValue "D4"  "5.339"
Set-TextValue "D5"  "0.05610"
Set-TextValue "D6"  "3.408"
Set-TextValue "D7"  "6.390"
Set-TextValue "D8"  "0.8159"
Set-TextValue "D9"  "0.9518"
Set-TextValue "D10" "0.1416"
Set-TextValue "D11" "0.07534"
Set-TextValue "D12" "0.03178"
Set-TextValue "D13" "0.03096"
Set-TextValue "D15" "3.558"
Set-TextValue "D16" "0.001606"
Set-TextValue "D17" "0.04703"

# --- Rows 18-24: coin/link/price/volume roll over one slot ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.006252"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.005086"
$ws.Range("E19").Value = "18HotbitTokenHTB"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.001032"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.0001499"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.748"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D23" "2.141"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D24" "0.01146"
$ws.Range("E24").Value = "23OneONEBestin24h"

# --- Remaining price-only refreshes (rows 25-50) ---
Set-TextValue "D25" "0.3255"
Set-TextValue "D26" "0.1288"
Set-TextValue "D28" "0.0002999"
Set-TextValue "D41" "0.006992"
Set-TextValue "D42" "0.1065"
Set-TextValue "D43" "0.003399"
Set-TextValue "D44" "0.008773"
Set-TextValue "D45" "0.00005726"

Set-TextValue "D47" "0.0005500"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"

Set-TextValue "D48" "0.7797"

Set-TextValue "D49" "0.1723"
$ws.Range("E49").Value = "48BOLOBOLO"

Set-TextValue "D50" "0.00002099"

Write-Host "Applied cryptos.xlsx price/symbol updates"
